$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.1
$ws.Range("O2").Value = 1.44
$ws.Range("P2").Value = 2.75
$ws.Range("AO2").Value = 8.5
$ws.Range("AW2").Value = 7.5

# Row 3
$ws.Range("M3").Value = 1.07
$ws.Range("O3").Value = 1.4

# Row 4
$ws.Range("G4").Value = 3.75
$ws.Range("I4").Value = 2.1
$ws.Range("J4").Value = 4.5
$ws.Range("L4").Value = 2.88
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
$ws.Range("O4").Value = 1.44
$ws.Range("P4").Value = 2.63
$ws.Range("Y4").Value = 15
$ws.Range("AI4").Value = 8.5
$ws.Range("AO4").Value = 23
$ws.Range("AX4").Value = 12

# Row 5
$ws.Range("Q5").Value = 2.5
$ws.Range("R5").Value = 1.5
$ws.Range("U5").Value = 2.62
$ws.Range("V5").Value = 1.41

# Row 6
$ws.Range("G6").Value = 5
$ws.Range("H6").Value = 3.3
$ws.Range("I6").Value = 1.8
$ws.Range("L6").Value = 2.5
$ws.Range("N6").Value = 7.5
$ws.Range("AI6").Value = 7.5
$ws.Range("AX6").Value = 10
$ws.Range("BA6").Value = 67

# Row 9
$ws.Range("M9").Value = 1.11
$ws.Range("N9").Value = 6.5
$ws.Range("Q9").Value = 2.6
$ws.Range("R9").Value = 1.48
$ws.Range("V9").Value = 1.67

# Row 11
$ws.Range("M11").Value = 1.08
$ws.Range("O11").Value = 1.5
$ws.Range("P11").Value = 2.37

# Row 12
$ws.Range("M12").Value = 1.08
$ws.Range("N12").Value = 8
$ws.Range("U12").Value = 2.37
$ws.Range("V12").Value = 1.5

# Row 13
$ws.Range("U13").Value = 1.69
